$d = $word.ActiveDocument

# 1. Bold the "Make search be just an icon" paragraph (text + paragraph mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Make search be just an icon*") {
        $p.Range.Bold = $true
        break
    }
}

# 2. Bold the "Take out the left and right container stuff" paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Take out the left and right container stuff*") {
        $p.Range.Bold = $true
        break
    }
}

# 3. The document's "_GoBack" bookmark moves from the "Make the buttons
#    purple..." paragraph to the end of the newly-added navbar paragraph.
#    Remove it from its old location first.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 4. Insert a new paragraph right after "Once it crosses the threshold..."
#    with the new task text, list level 3 / numId 1, and carry the
#    "_GoBack" bookmark at the end of its run (point bookmark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Once it crosses the threshold*") {
        $insertPos = $p.Range.End - 1
        $insertRange = $d.Range($insertPos, $insertPos)
        $newParagraphXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="3"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Fix the layout and the structure of the elements within the navbar</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
        $insertRange.InsertXML($newParagraphXml)
        break
    }
}
